$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was collected; insert it as a new row at 161,
# pushing the existing rows (161-197) down to (162-198).
$ws.Rows.Item(161).Insert()

$ws.Cells.Item(161, 1).Value = 4
$ws.Cells.Item(161, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(161, 3).Value = "Los Lagos"
$ws.Cells.Item(161, 4).Value = 44476
$ws.Cells.Item(161, 5).Value = 10
$ws.Cells.Item(161, 6).Value = 100112023
$ws.Cells.Item(161, 7).Value = "Brócoli"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 500
$ws.Cells.Item(161, 11).Value = 1100
$ws.Cells.Item(161, 12).Value = 1200
$ws.Cells.Item(161, 13).Value = 1150
$ws.Cells.Item(161, 14).Value = "$/unidad"
$ws.Cells.Item(161, 15).Value = "Región Metropolitana"
$ws.Cells.Item(161, 16).Value = 1150
$ws.Cells.Item(161, 17).Value = 1
$ws.Cells.Item(161, 18).Value = "Hortaliza"
